# Update a couple of unit_cost values on the "endofpipe" sheet, then move the
# user's active selection: inputdisp's selection moves from G3 to B3 (and is
# no longer the active tab), while endofpipe becomes the active tab with its
# selection at G2.

$wb = $excel.ActiveWorkbook

$wsDisp = $wb.Worksheets.Item("inputdisp")
$wsPipe = $wb.Worksheets.Item("endofpipe")

# Data edits: endofpipe!G2 10 -> 6, endofpipe!G3 6.5 -> 2
$wsPipe.Range("G2").Value = 6
$wsPipe.Range("G3").Value = 2

# Move inputdisp's selection first (while it is still active) so the sheet
# activation order ends with "endofpipe" as the final active tab.
$wsDisp.Activate()
$wsDisp.Range("B3").Select()

# endofpipe becomes the active sheet/tab, selection parked at G2.
$wsPipe.Activate()
$wsPipe.Range("G2").Select()
